$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("E9").Value = 24

# Row 17
$ws.Range("E17").Value = 111
$ws.Range("F17").Value = 50
$ws.Range("H17").Value = 50

# Row 35
$ws.Range("E35").Value = 7

# Row 36
$ws.Range("E36").Value = 94
$ws.Range("F36").Value = 42
$ws.Range("H36").Value = 42

# Row 37
$ws.Range("E37").Value = 49

# Row 38
$ws.Range("E38").Value = 69

# Row 39
$ws.Range("E39").Value = 25

# Row 44
$ws.Range("E44").Value = 28

# Row 45
$ws.Range("E45").Value = 24

# Row 49
$ws.Range("E49").Value = 63

# Row 63
$ws.Range("E63").Value = 33
$ws.Range("F63").Value = 11
$ws.Range("H63").Value = 11

# Row 71
$ws.Range("E71").Value = 32

# Row 75
$ws.Range("E75").Value = 13
$ws.Range("F75").Value = 6
$ws.Range("H75").Value = 6

# Row 79
$ws.Range("E79").Value = 36

# Row 80
$ws.Range("E80").Value = 25

# Row 89
$ws.Range("E89").Value = 33

$wb.Save()
